$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-08 Sunday", 2) | Out-Null
$d.Content.Find.Execute("975÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "317÷8=", 2) | Out-Null
$d.Content.Find.Execute("743÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "269÷2=", 2) | Out-Null
$d.Content.Find.Execute("829÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "710÷5=", 2) | Out-Null
$d.Content.Find.Execute("776÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "253÷9=", 2) | Out-Null
$d.Content.Find.Execute("520÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "886÷3=", 2) | Out-Null
$d.Content.Find.Execute("239÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "815÷5=", 2) | Out-Null
$d.Content.Find.Execute("721÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "513÷8=", 2) | Out-Null
$d.Content.Find.Execute("985÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "707÷9=", 2) | Out-Null
$d.Content.Find.Execute("715÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "841÷8=", 2) | Out-Null
$d.Content.Find.Execute("293÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "603÷2=", 2) | Out-Null
$d.Content.Find.Execute("856÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "156÷9=", 2) | Out-Null
$d.Content.Find.Execute("471÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "929÷3=", 2) | Out-Null
$d.Content.Find.Execute("367÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "809÷3=", 2) | Out-Null
$d.Content.Find.Execute("280÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "231÷3=", 2) | Out-Null
$d.Content.Find.Execute("891÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "468÷6=", 2) | Out-Null
$d.Content.Find.Execute("197÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "101÷2=", 2) | Out-Null
$d.Content.Find.Execute("232÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "517÷6=", 2) | Out-Null
$d.Content.Find.Execute("926÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "465÷4=", 2) | Out-Null
$d.Content.Find.Execute("546÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "648÷3=", 2) | Out-Null
$d.Content.Find.Execute("784÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "530÷8=", 2) | Out-Null
$d.Content.Find.Execute("487÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "304÷8=", 2) | Out-Null
$d.Content.Find.Execute("641÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "765÷8=", 2) | Out-Null
$d.Content.Find.Execute("692÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "994÷9=", 2) | Out-Null
$d.Content.Find.Execute("451÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "467÷4=", 2) | Out-Null
$d.Content.Find.Execute("133÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "786÷4=", 2) | Out-Null
